# Backup QR Scanner data - 26/07/2025, 10:36:25 PM
$wb = $excel.ActiveWorkbook

# Remove the extra "Anatomy_2" and "Anatomy_3" sheets, keeping only "Anatomy".
foreach ($sheetName in @("Anatomy_2", "Anatomy_3")) {
    $sheet = $wb.Worksheets.Item($sheetName)
    if ($sheet) {
        $sheet.Delete()
    }
}

# Update the remaining log entry on the "Anatomy" sheet.
$ws = $wb.Worksheets.Item("Anatomy")

# A2 holds a numeric-looking ID that must stay text (matches original "str" cell type).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "424333"
$ws.Range("A2").Style = "Normal"

$ws.Range("C2").Value = "26/07/2025"
$ws.Range("D2").Value = "22:35:50"
